$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values (e.g. "8.00", "568.18")
# are not auto-converted to numbers by Excel's type inference when set via .Value
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.555.06"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "2.410.34"
$ws.Range("E3").Value = "  +0.14%  "

$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").Value = "568.18"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").Value = "143.49"
$ws.Range("E6").Value = "  +4.02%  "

$ws.Range("E7").Value = "  -0.51%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "2.421.30"
$ws.Range("E9").Value = "  +1.35%  "

$ws.Range("D10").Value = "0.108"
$ws.Range("E10").Value = "  +2.82%  "

$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  +0.52%  "

$ws.Range("D12").Value = "5.22"
$ws.Range("E12").Value = "  +3.84%  "

$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  +4.16%  "

$ws.Range("D14").Value = "26.50"
$ws.Range("E14").Value = "  +3.08%  "

$ws.Range("E15").Value = "  +3.55%  "

$ws.Range("D16").Value = "2.842.98"
$ws.Range("E16").Value = "  -0.18%  "

$ws.Range("D17").Value = "61.442.61"
$ws.Range("E17").Value = "  +1.10%  "

$ws.Range("D18").Value = "2.421.07"
$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  -0.04%  "

$ws.Range("D20").Value = "10.66"
$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").Value = "324.19"
$ws.Range("E21").Value = "  +0.80%  "

$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("E23").Value = "  -0.99%  "

$ws.Range("D24").Value = "1.98"
$ws.Range("E24").Value = "  +9.47%  "

$ws.Range("E25").Value = "  -0.15%  "

$ws.Range("D26").Value = "65.00"
$ws.Range("E26").Value = "  +1.25%  "

$ws.Range("D27").Value = "615.65"
$ws.Range("E27").Value = "  +8.67%  "

$ws.Range("D28").Value = "8.29"
$ws.Range("E28").Value = "  +1.30%  "

$ws.Range("E29").Value = "  +4.39%  "

$ws.Range("D31").Value = "8.05"
$ws.Range("E31").Value = "  +1.93%  "

$ws.Range("E32").Value = "  +5.60%  "

$ws.Range("D33").Value = "1.81"
$ws.Range("E33").Value = "  +0.53%  "

$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("E35").Value = "  +5.50%  "

$ws.Range("D36").Value = "0.996"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "0.373"
$ws.Range("E37").Value = "  +1.77%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").Value = "4.63"
$ws.Range("E38").Value = "  +2.85%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").Value = "151.86"
$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("D40").Value = "5.37"
$ws.Range("E40").Value = "  +6.05%  "

$ws.Range("D41").Value = "18.40"
$ws.Range("E41").Value = "  +1.65%  "

$ws.Range("E42").Value = "  +9.53%  "

$ws.Range("E43").Value = "  +4.42%  "

$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "41.88"
$ws.Range("E45").Value = "  +1.86%  "

$ws.Range("D46").Value = "0.0₆0282"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("D47").Value = "142.72"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").Value = "3.55"
$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("D49").Value = "19.94"
$ws.Range("E49").Value = "  +4.93%  "

$ws.Range("E50").Value = "  +2.20%  "

$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").Value = "  +2.84%  "

# Restore column D to the default (General) style so no residual number formatting remains
$ws.Range("D2:D51").Style = "Normal"
